# Update "想去人数" (want-to-go count) figures that changed between scrapes.
# Affects two sheets that both list the exhibition ("展览") events:
#   - "展览"    (the exhibition-only sheet)
#   - "全部类型" (the combined/all-types sheet)

$wb = $excel.ActiveWorkbook

$sheetExpo = $wb.Worksheets.Item("展览")
$sheetAll  = $wb.Worksheets.Item("全部类型")

# Row -> new value updates for the "展览" sheet
$sheetExpo.Range("F2").Value  = 804
$sheetExpo.Range("F5").Value  = 1071
$sheetExpo.Range("F8").Value  = 210
$sheetExpo.Range("F9").Value  = 386
$sheetExpo.Range("F15").Value = 12505
$sheetExpo.Range("F16").Value = 141
$sheetExpo.Range("F17").Value = 5498

# Same events, different row numbers, on the "全部类型" sheet
$sheetAll.Range("F2").Value  = 804
$sheetAll.Range("F7").Value  = 1071
$sheetAll.Range("F10").Value = 210
$sheetAll.Range("F11").Value = 386
$sheetAll.Range("F17").Value = 12505
$sheetAll.Range("F20").Value = 141
$sheetAll.Range("F21").Value = 5498
